# Apply "added fresh data for registration" changes to TestData.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegDetails")

# Update the registration rows with fresh names / emails.
# Column A (FirstName) and Column E (Password) are unchanged.
# Columns B (MiddleName/Initial), C (LastName) and D (EmailAddress) get new values.
$ws.Range("B2").Value = "Ronald"
$ws.Range("C2").Value = "Delver"
$ws.Range("D2").Value = "a25@email.com"

$ws.Range("B3").Value = "Peter"
$ws.Range("C3").Value = "Con"
$ws.Range("D3").Value = "a26@email.com"

$ws.Range("B4").Value = "Asult"
$ws.Range("C4").Value = "Bolswa"
$ws.Range("D4").Value = "a27@email.com"

$ws.Range("B5").Value = "Jimmy"
$ws.Range("C5").Value = "Lever"
$ws.Range("D5").Value = "a28@email.com"

# Update the selected cell in the worksheet view from E3 to E4.
$ws.Range("E4").Select()

# Update the workbook window size recorded in the workbook view.
$win = $excel.ActiveWindow
$win.Width = 19380
$win.Height = 5955
$excel.Width = 19380
$excel.Height = 5955
